# Add a new "Correction " column to the Card17 sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card17")

# Copy the style of the existing "Event " header (M1, bold/centered/bordered)
# onto the new N1 header cell before setting its text.
$ws.Range("M1").Copy() | Out-Null
$ws.Range("N1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Fix header text of column M ("Event " -> "Event") and set the new header text.
$ws.Range("M1").Value = "Event"
$ws.Range("N1").Value = "Correction "

# Fill in the data rows: M column gets "nan" text where it was previously blank,
# and the new N column gets blank placeholder cells (materialized, but empty)
# for each data row, matching the pattern used on sibling "Card" sheets.
for ($r = 2; $r -le 12; $r++) {
    $ws.Cells.Item($r, 13).Value = "nan"    # column M
    $ws.Cells.Item($r, 14).Style = "Normal" # column N (materialize blank cell)
}
